$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text, preventing Excel from
# reinterpreting values such as "290.87" as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '22.371.55'
$ws.Range('E2').Value = '  -4.73%  '
$ws.Range('D3').Value = '1.573.11'
$ws.Range('E3').Value = '  -4.54%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').Value = '290.87'
$ws.Range('E6').Value = '  -2.92%  '
$ws.Range('D7').Value = '0.3777'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '49.82'
$ws.Range('E8').Value = '  -2.31%  '
$ws.Range('D9').Value = '0.3436'
$ws.Range('E9').Value = '  -1.92%  '
$ws.Range('D10').Value = '1.171'
$ws.Range('E10').Value = '  -4.31%  '
$ws.Range('D11').Value = '0.07693'
$ws.Range('E11').Value = '  -4.66%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = '21.44'
$ws.Range('E13').Value = '  -2.32%  '
$ws.Range('D14').Value = '6.033'
$ws.Range('E14').Value = '  -4.75%  '
$ws.Range('D15').Value = '6.979'
$ws.Range('E15').Value = '  -3.93%  '
$ws.Range('D16').Value = '0.00001143'
$ws.Range('E16').Value = '  -5.07%  '
$ws.Range('D17').Value = '1.576.25'
$ws.Range('E17').Value = '  -3.96%  '
$ws.Range('D18').Value = '90.68'
$ws.Range('E18').Value = '  -4.83%  '
$ws.Range('D19').Value = '0.06742'
$ws.Range('E19').Value = '  -3.47%  '
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '6.288'
$ws.Range('E21').Value = '  -4.92%  '
$ws.Range('D22').Value = '16.74'
$ws.Range('E22').Value = '  -4.10%  '
$ws.Range('D23').Value = '0.5311'
$ws.Range('E23').Value = '  -8.24%  '
$ws.Range('D24').Value = '12.04'
$ws.Range('E24').Value = '  -3.12%  '
$ws.Range('D25').Value = '22.374.60'
$ws.Range('E25').Value = '  -4.66%  '
$ws.Range('D26').Value = '2.397'
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').Value = '2.805'
$ws.Range('E27').Value = '  -7.38%  '
$ws.Range('D28').Value = '20.36'
$ws.Range('E28').Value = '  -3.35%  '
$ws.Range('D29').Value = '145.21'
$ws.Range('E29').Value = '  -4.14%  '
$ws.Range('D30').Value = '5.069'
$ws.Range('E30').Value = '  -2.30%  '
$ws.Range('D31').Value = '126.75'
$ws.Range('E31').Value = '  -3.85%  '
$ws.Range('D32').Value = '1.748.29'
$ws.Range('E32').Value = '  -4.57%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '1.034'
$ws.Range('E33').Value = '  +5.29%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '6.293'
$ws.Range('E34').Value = '  -9.14%  '
$ws.Range('D35').Value = '2.028'
$ws.Range('E35').Value = '  -5.60%  '
$ws.Range('D36').Value = '10.13'
$ws.Range('E36').Value = '  -6.62%  '
$ws.Range('D37').Value = '0.08597'
$ws.Range('E37').Value = '  -2.21%  '
$ws.Range('D38').Value = '0.02578'
$ws.Range('E38').Value = '  -5.59%  '
$ws.Range('D39').Value = '0.2333'
$ws.Range('E39').Value = '  -3.45%  '
$ws.Range('D40').Value = '5.600'
$ws.Range('E40').Value = '  -5.02%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '0.06513'
$ws.Range('E41').Value = '  -4.85%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '1.319'
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('D43').Value = '11.80'
$ws.Range('E43').Value = '  -7.68%  '
$ws.Range('D44').Value = '0.6482'
$ws.Range('E44').Value = '  -5.85%  '
$ws.Range('D45').Value = '14.32'
$ws.Range('E45').Value = '  -7.47%  '
$ws.Range('D46').Value = '0.9998'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').Value = '0.6066'
$ws.Range('E47').Value = '  -4.73%  '
$ws.Range('D48').Value = '3.774'
$ws.Range('E48').Value = '  -3.67%  '
$ws.Range('D49').Value = '1.319'
$ws.Range('E49').Value = '  +6.28%  '
$ws.Range('D50').Value = '2.111'
$ws.Range('E50').Value = '  -6.52%  '
$ws.Range('D51').Value = '124.95'
$ws.Range('E51').Value = '  -1.83%  '
